$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-64 down to 17-65
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new weekly price-report record
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44497
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100108
$ws.Range("H16").Value = "Tropicales y subtropicales"
$ws.Range("I16").Value = 100108002
$ws.Range("J16").Value = "Mango"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 30
$ws.Range("N16").Value = 8000
$ws.Range("O16").Value = 8500
$ws.Range("P16").Value = 8250
$ws.Range("Q16").Value = "$/bandeja 4 kilos"
$ws.Range("R16").Value = "Perú"
$ws.Range("S16").Value = 2062
$ws.Range("T16").Value = 4

Write-Output "row inserted"
